$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 23, shifting the existing rows 23-41 down to 24-42.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new price-report record.
$ws.Cells.Item(23, 1).Value = 5
$ws.Cells.Item(23, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(23, 3).Value = "Maule"
$ws.Cells.Item(23, 4).Value = 44539
$ws.Cells.Item(23, 5).Value = 7
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100101
$ws.Cells.Item(23, 8).Value = "Berries"
$ws.Cells.Item(23, 9).Value = 100101001
$ws.Cells.Item(23, 10).Value = "Arándano (blue)"
$ws.Cells.Item(23, 11).Value = "Sin especificar"
$ws.Cells.Item(23, 12).Value = "Segunda"
$ws.Cells.Item(23, 13).Value = 150
$ws.Cells.Item(23, 14).Value = 3000
$ws.Cells.Item(23, 15).Value = 3000
$ws.Cells.Item(23, 16).Value = 3000
$ws.Cells.Item(23, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(23, 18).Value = "Provincia de Linares"
$ws.Cells.Item(23, 19).Value = 1500
$ws.Cells.Item(23, 20).Value = 2
